# Scheduled-runner style refresh of market-price-derived columns
# (currentAveragePrice / currentAveragePriceNQ / currentAveragePriceHQ /
#  LevePriceNQ / LevePriceHQ / LeveProfitNQ / LeveProfitHQ) across the
# ALC / ARM / BSM / CRP / CUL / GSM / LTW / WVR sheets.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H55").Value = 2245.7
$ws.Range("I55").Value = 179.71428
$ws.Range("J55").Value = 7066.3335
$ws.Range("K55").Value = 179.71428
$ws.Range("L55").Value = 7066.3335
$ws.Range("M55").Value = 34.28572
$ws.Range("N55").Value = -7494.3335

$ws.Range("H137").Value = 2979.5454
$ws.Range("J137").Value = 5015.933
$ws.Range("L137").Value = 15047.799
$ws.Range("N137").Value = -20147.799

$ws.Range("H138").Value = 2632.6667
$ws.Range("I138").Value = 940.2222
$ws.Range("J138").Value = 3140.4
$ws.Range("K138").Value = 2820.6666
$ws.Range("L138").Value = 9421.200000000001
$ws.Range("M138").Value = 2319.3334
$ws.Range("N138").Value = -19701.2


$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H122").Value = 1401.8334
$ws.Range("I122").Value = 1137.3334
$ws.Range("K122").Value = 3412.0002
$ws.Range("M122").Value = -962.0001999999999


$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H76").Value = 80001
$ws.Range("J76").Value = 80001
$ws.Range("L76").Value = 80001
$ws.Range("N76").Value = -80631

$ws.Range("H79").Value = 80001
$ws.Range("J79").Value = 80001
$ws.Range("L79").Value = 80001
$ws.Range("N79").Value = -82185

$ws.Range("H82").Value = 54416.668
$ws.Range("J82").Value = 100000
$ws.Range("L82").Value = 100000
$ws.Range("N82").Value = -100766

$ws.Range("H85").Value = 54416.668
$ws.Range("J85").Value = 100000
$ws.Range("L85").Value = 100000
$ws.Range("N85").Value = -102652

$ws.Range("H94").Value = 896
$ws.Range("I94").Value = 896
$ws.Range("K94").Value = 896
$ws.Range("M94").Value = -445

$ws.Range("H96").Value = 34893.727
$ws.Range("J96").Value = 76108.25
$ws.Range("L96").Value = 76108.25
$ws.Range("N96").Value = -81600.25


$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 3512.6
$ws.Range("I7").Value = 114.57143
$ws.Range("J7").Value = 11441.333
$ws.Range("K7").Value = 114.57143
$ws.Range("L7").Value = 11441.333
$ws.Range("M7").Value = -1.571430000000007
$ws.Range("N7").Value = -11667.333

$ws.Range("H92").Value = 42996.5
$ws.Range("J92").Value = 42996.5
$ws.Range("L92").Value = 42996.5
$ws.Range("N92").Value = -47988.5

$ws.Range("H127").Value = 119000
$ws.Range("J127").Value = 119000
$ws.Range("L127").Value = 119000
$ws.Range("N127").Value = -128920


$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H3").Value = 5753.625
$ws.Range("I3").Value = 1006.44446
$ws.Range("J3").Value = 11857.143
$ws.Range("K3").Value = 3019.33338
$ws.Range("L3").Value = 35571.429
$ws.Range("M3").Value = -2907.33338
$ws.Range("N3").Value = -35795.429

$ws.Range("H5").Value = 903.5833
$ws.Range("I5").Value = 849.36365
$ws.Range("J5").Value = 1500
$ws.Range("K5").Value = 2548.09095
$ws.Range("L5").Value = 4500
$ws.Range("M5").Value = -2436.09095
$ws.Range("N5").Value = -4724

$ws.Range("H11").Value = 2088.0833
$ws.Range("J11").Value = 840
$ws.Range("L11").Value = 2520
$ws.Range("N11").Value = -2800

$ws.Range("H107").Value = 394.42856
$ws.Range("J107").Value = 455.4
$ws.Range("L107").Value = 1366.2
$ws.Range("N107").Value = -5206.2

$ws.Range("H132").Value = 1934.1428
$ws.Range("I132").Value = 2025.9
$ws.Range("J132").Value = 1850.7273
$ws.Range("K132").Value = 18233.1
$ws.Range("L132").Value = 16656.5457
$ws.Range("M132").Value = -15703.1
$ws.Range("N132").Value = -21716.5457

$ws.Range("H134").Value = 4868.8823
$ws.Range("I134").Value = 3165.1667
$ws.Range("K134").Value = 9495.500100000001
$ws.Range("M134").Value = -4425.500100000001

$ws.Range("H135").Value = 903.5833
$ws.Range("I135").Value = 849.36365
$ws.Range("J135").Value = 1500
$ws.Range("K135").Value = 7644.27285
$ws.Range("L135").Value = 13500
$ws.Range("M135").Value = -5109.27285
$ws.Range("N135").Value = -18570


$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H7").Value = 62502.25
$ws.Range("I7").Value = 10000
$ws.Range("K7").Value = 10000
$ws.Range("M7").Value = -9888

$ws.Range("H8").Value = 62502.25
$ws.Range("I8").Value = 10000
$ws.Range("K8").Value = 10000
$ws.Range("M8").Value = -9861

$ws.Range("H22").Value = 10009
$ws.Range("I22").Value = 0
$ws.Range("J22").Value = 10009
$ws.Range("K22").Value = 0
$ws.Range("L22").Value = 10009
$ws.Range("M22").ClearContents()
$ws.Range("N22").Value = -11067

$ws.Range("H62").Value = 100437
$ws.Range("J62").Value = 100437
$ws.Range("L62").Value = 100437
$ws.Range("N62").Value = -101809

$ws.Range("H65").Value = 100437
$ws.Range("J65").Value = 100437
$ws.Range("L65").Value = 301311
$ws.Range("N65").Value = -308175

$ws.Range("H70").Value = 4940
$ws.Range("I70").Value = 4925
$ws.Range("K70").Value = 4925
$ws.Range("M70").Value = -4655

$ws.Range("H73").Value = 4940
$ws.Range("I73").Value = 4925
$ws.Range("K73").Value = 4925
$ws.Range("M73").Value = -3989

$ws.Range("H92").Value = 20784.285
$ws.Range("J92").Value = 20784.285
$ws.Range("L92").Value = 20784.285
$ws.Range("N92").Value = -24528.285

$ws.Range("H107").Value = 1528.4
$ws.Range("I107").Value = 1420.4445
$ws.Range("J107").Value = 2500
$ws.Range("K107").Value = 1420.4445
$ws.Range("L107").Value = 2500
$ws.Range("M107").Value = 499.5554999999999
$ws.Range("N107").Value = -6340


$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H5").Value = 0
$ws.Range("I5").Value = 0
$ws.Range("J5").Value = 0
$ws.Range("K5").Value = 0
$ws.Range("L5").Value = 0
$ws.Range("M5").ClearContents()
$ws.Range("N5").ClearContents()

$ws.Range("H20").Value = 0
$ws.Range("I20").Value = 0
$ws.Range("J20").Value = 0
$ws.Range("K20").Value = 0
$ws.Range("L20").Value = 0
$ws.Range("M20").ClearContents()
$ws.Range("N20").ClearContents()

$ws.Range("H21").Value = 10006.533

$ws.Range("H24").Value = 60006.5
$ws.Range("I24").Value = 40006
$ws.Range("J24").Value = 80007
$ws.Range("K24").Value = 40006
$ws.Range("L24").Value = 80007
$ws.Range("M24").Value = -39663
$ws.Range("N24").Value = -80693

$ws.Range("H93").Value = 100002080
$ws.Range("I93").Value = 111112980
$ws.Range("K93").Value = 111112980
$ws.Range("M93").Value = -111111732

$ws.Range("H101").Value = 51551.57
$ws.Range("J101").Value = 51551.57
$ws.Range("L101").Value = 51551.57
$ws.Range("N101").Value = -58041.57


$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H24").Value = 0
$ws.Range("J24").Value = 0
$ws.Range("L24").Value = 0
$ws.Range("N24").ClearContents()

$ws.Range("H28").Value = 0
$ws.Range("I28").Value = 0
$ws.Range("J28").Value = 0
$ws.Range("K28").Value = 0
$ws.Range("L28").Value = 0
$ws.Range("M28").ClearContents()
$ws.Range("N28").ClearContents()

$ws.Range("H30").Value = 20009
$ws.Range("I30").Value = 20009
$ws.Range("K30").Value = 20009
$ws.Range("M30").Value = -19902

$ws.Range("H64").Value = 64997
$ws.Range("J64").Value = 64997
$ws.Range("L64").Value = 64997
$ws.Range("N64").Value = -65493

$ws.Range("H67").Value = 64997
$ws.Range("J67").Value = 64997
$ws.Range("L67").Value = 64997
$ws.Range("N67").Value = -66713

$ws.Range("H81").Value = 1199.75
$ws.Range("I81").Value = 1199.75
$ws.Range("K81").Value = 2399.5
$ws.Range("M81").Value = -1338.5

$ws.Range("H84").Value = 1199.75
$ws.Range("I84").Value = 1199.75
$ws.Range("K84").Value = 11997.5
$ws.Range("M84").Value = -6693.5

$ws.Range("H93").Value = 77882
$ws.Range("I93").Value = 0
$ws.Range("J93").Value = 77882
$ws.Range("K93").Value = 0
$ws.Range("L93").Value = 77882
$ws.Range("M93").ClearContents()
$ws.Range("N93").Value = -82874

$ws.Range("H96").Value = 3749.75

$ws.Range("H103").Value = 52782.145
$ws.Range("J103").Value = 52782.145
$ws.Range("L103").Value = 52782.145
$ws.Range("N103").Value = -55126.145

